$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "70.059.58"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "3.770.22"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "624.20"
$ws.Range("E5").Value = "  +0.80%  "
Set-TextValue "D6" "180.18"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").Value = "3.767.43"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  -5.03%  "
$ws.Range("E12").Value = "  -2.67%  "
Set-TextValue "D13" "41.20"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "4.387.34"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").Value = "3.766.74"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "70.114.49"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("E18").Value = "  -0.26%  "
Set-TextValue "D19" "7.63"
$ws.Range("E19").Value = "  +0.92%  "
Set-TextValue "D20" "16.77"
$ws.Range("E20").Value = "  -0.89%  "
Set-TextValue "D21" "507.41"
$ws.Range("E21").Value = "  -2.73%  "
Set-TextValue "D22" "9.50"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("E24").Value = "  -1.35%  "
Set-TextValue "D25" "87.21"
$ws.Range("E25").Value = "  -1.80%  "
Set-TextValue "D26" "13.20"
$ws.Range("E26").Value = "  -1.87%  "
Set-TextValue "D27" "11.20"
$ws.Range("E27").Value = "  +0.69%  "
Set-TextValue "D28" "0.0000139"
$ws.Range("E28").Value = "  +26.61%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -2.25%  "
Set-TextValue "D31" "2.96"
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("E32").Value = "  -3.50%  "
Set-TextValue "D33" "31.42"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +0.20%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +4.85%  "
$ws.Range("E37").Value = "  +1.56%  "
Set-TextValue "D38" "0.334"
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  -3.63%  "
Set-TextValue "D41" "50.56"
$ws.Range("E41").Value = "  -1.65%  "
Set-TextValue "D42" "45.12"
$ws.Range("E42").Value = "  -2.08%  "
Set-TextValue "D43" "424.27"
$ws.Range("E43").Value = "  -0.83%  "
Set-TextValue "D44" "8.74"
$ws.Range("E44").Value = "  -0.99%  "
Set-TextValue "D45" "2.85"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "3.012.43"
$ws.Range("E46").Value = "  -3.25%  "
Set-TextValue "D47" "0.0366"
$ws.Range("E47").Value = "  -0.98%  "
Set-TextValue "D48" "27.40"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D49" "1.00"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "138.38"
$ws.Range("E50").Value = "  -1.60%  "
Set-TextValue "D51" "2.53"
$ws.Range("E51").Value = "  +1.91%  "
